$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Row = 8; Col = "A"; Value = 46029 },
    @{ Row = 8; Col = "D"; Value = 152.5 },
    @{ Row = 8; Col = "E"; Value = 149.77000000000001 },
    @{ Row = 8; Col = "F"; Value = 159.78 },
    @{ Row = 8; Col = "G"; Value = 149.79 },
    @{ Row = 9; Col = "A"; Value = 46029 },
    @{ Row = 9; Col = "D"; Value = 152.5 },
    @{ Row = 9; Col = "E"; Value = 149.77000000000001 },
    @{ Row = 9; Col = "F"; Value = 159.78 },
    @{ Row = 9; Col = "G"; Value = 149.79 },
    @{ Row = 10; Col = "A"; Value = 46029 },
    @{ Row = 10; Col = "D"; Value = 154.68 },
    @{ Row = 10; Col = "E"; Value = 151.57 },
    @{ Row = 10; Col = "F"; Value = 161.57 },
    @{ Row = 10; Col = "G"; Value = 151.97 },
    @{ Row = 11; Col = "A"; Value = 46028 },
    @{ Row = 11; Col = "D"; Value = 152.03 },
    @{ Row = 11; Col = "E"; Value = 149.57 },
    @{ Row = 11; Col = "F"; Value = 159.57 },
    @{ Row = 11; Col = "G"; Value = 149.59 },
    @{ Row = 12; Col = "A"; Value = 46028 },
    @{ Row = 12; Col = "D"; Value = 152.03 },
    @{ Row = 12; Col = "E"; Value = 149.57 },
    @{ Row = 12; Col = "F"; Value = 159.57 },
    @{ Row = 12; Col = "G"; Value = 149.59 },
    @{ Row = 13; Col = "A"; Value = 46028 },
    @{ Row = 13; Col = "D"; Value = 154.24 },
    @{ Row = 13; Col = "E"; Value = 151.4 },
    @{ Row = 13; Col = "F"; Value = 161.4 },
    @{ Row = 13; Col = "G"; Value = 151.80000000000001 },
    @{ Row = 17; Col = "A"; Value = 46029 },
    @{ Row = 17; Col = "D"; Value = 158 },
    @{ Row = 17; Col = "E"; Value = 154.63 },
    @{ Row = 17; Col = "F"; Value = 164.63 },
    @{ Row = 18; Col = "A"; Value = 46028 },
    @{ Row = 18; Col = "D"; Value = 158.13999999999999 },
    @{ Row = 18; Col = "E"; Value = 155.03 },
    @{ Row = 18; Col = "F"; Value = 165.03 },
    @{ Row = 22; Col = "A"; Value = 46029 },
    @{ Row = 22; Col = "D"; Value = 154.02000000000001 },
    @{ Row = 22; Col = "E"; Value = 151.02000000000001 },
    @{ Row = 22; Col = "F"; Value = 160.62 },
    @{ Row = 22; Col = "G"; Value = 152.09 },
    @{ Row = 23; Col = "A"; Value = 46029 },
    @{ Row = 23; Col = "D"; Value = 158.76 },
    @{ Row = 23; Col = "E"; Value = 156.47999999999999 },
    @{ Row = 23; Col = "F"; Value = 166.48 },
    @{ Row = 24; Col = "A"; Value = 46029 },
    @{ Row = 24; Col = "D"; Value = 158.93 },
    @{ Row = 24; Col = "E"; Value = 157.09 },
    @{ Row = 24; Col = "F"; Value = 167.09 },
    @{ Row = 25; Col = "A"; Value = 46029 },
    @{ Row = 25; Col = "D"; Value = 158.91 },
    @{ Row = 25; Col = "E"; Value = 156.61000000000001 },
    @{ Row = 25; Col = "F"; Value = 166.61 },
    @{ Row = 25; Col = "G"; Value = 156.72999999999999 },
    @{ Row = 26; Col = "A"; Value = 46029 },
    @{ Row = 26; Col = "D"; Value = 158.56 },
    @{ Row = 26; Col = "E"; Value = 158.22 },
    @{ Row = 26; Col = "F"; Value = 168.22 },
    @{ Row = 27; Col = "A"; Value = 46028 },
    @{ Row = 27; Col = "D"; Value = 153.77000000000001 },
    @{ Row = 27; Col = "E"; Value = 150.81 },
    @{ Row = 27; Col = "F"; Value = 160.41 },
    @{ Row = 27; Col = "G"; Value = 151.88 },
    @{ Row = 28; Col = "A"; Value = 46028 },
    @{ Row = 28; Col = "D"; Value = 158.55000000000001 },
    @{ Row = 28; Col = "E"; Value = 156.31 },
    @{ Row = 28; Col = "F"; Value = 166.31 },
    @{ Row = 29; Col = "A"; Value = 46028 },
    @{ Row = 29; Col = "D"; Value = 158.71 },
    @{ Row = 29; Col = "E"; Value = 156.91 },
    @{ Row = 29; Col = "F"; Value = 166.91 },
    @{ Row = 30; Col = "A"; Value = 46028 },
    @{ Row = 30; Col = "D"; Value = 158.69999999999999 },
    @{ Row = 30; Col = "E"; Value = 156.41999999999999 },
    @{ Row = 30; Col = "F"; Value = 166.42 },
    @{ Row = 30; Col = "G"; Value = 156.55000000000001 },
    @{ Row = 31; Col = "A"; Value = 46028 },
    @{ Row = 31; Col = "D"; Value = 158.35 },
    @{ Row = 31; Col = "E"; Value = 158.04 },
    @{ Row = 31; Col = "F"; Value = 168.04 },
    @{ Row = 35; Col = "A"; Value = 46029 },
    @{ Row = 35; Col = "D"; Value = 151.69 },
    @{ Row = 35; Col = "E"; Value = 149.43 },
    @{ Row = 35; Col = "F"; Value = 158.43 },
    @{ Row = 36; Col = "A"; Value = 46028 },
    @{ Row = 36; Col = "D"; Value = 151.80000000000001 },
    @{ Row = 36; Col = "E"; Value = 149.81 },
    @{ Row = 36; Col = "F"; Value = 158.81 },
    @{ Row = 40; Col = "A"; Value = 46029 },
    @{ Row = 40; Col = "D"; Value = 159.13 },
    @{ Row = 40; Col = "E"; Value = 157.28 },
    @{ Row = 40; Col = "F"; Value = 167.28 },
    @{ Row = 41; Col = "A"; Value = 46029 },
    @{ Row = 41; Col = "D"; Value = 158.85 },
    @{ Row = 41; Col = "E"; Value = 157.69999999999999 },
    @{ Row = 41; Col = "F"; Value = 167.7 },
    @{ Row = 42; Col = "A"; Value = 46028 },
    @{ Row = 42; Col = "D"; Value = 159.22999999999999 },
    @{ Row = 42; Col = "E"; Value = 157.4 },
    @{ Row = 42; Col = "F"; Value = 167.4 },
    @{ Row = 43; Col = "A"; Value = 46028 },
    @{ Row = 43; Col = "D"; Value = 158.94999999999999 },
    @{ Row = 43; Col = "E"; Value = 157.82 },
    @{ Row = 43; Col = "F"; Value = 167.82 },
    @{ Row = 47; Col = "A"; Value = 46029 },
    @{ Row = 47; Col = "D"; Value = 152.53 },
    @{ Row = 47; Col = "E"; Value = 151.29 },
    @{ Row = 47; Col = "F"; Value = 161.29 },
    @{ Row = 48; Col = "A"; Value = 46029 },
    @{ Row = 48; Col = "D"; Value = 152.19999999999999 },
    @{ Row = 48; Col = "E"; Value = 151.24 },
    @{ Row = 48; Col = "F"; Value = 161.24 },
    @{ Row = 49; Col = "A"; Value = 46028 },
    @{ Row = 49; Col = "D"; Value = 152.71 },
    @{ Row = 49; Col = "E"; Value = 151.57 },
    @{ Row = 49; Col = "F"; Value = 161.57 },
    @{ Row = 50; Col = "A"; Value = 46028 },
    @{ Row = 50; Col = "D"; Value = 152.38 },
    @{ Row = 50; Col = "E"; Value = 151.52000000000001 },
    @{ Row = 50; Col = "F"; Value = 161.52000000000001 },
    @{ Row = 54; Col = "A"; Value = 46029 },
    @{ Row = 54; Col = "D"; Value = 168.43 },
    @{ Row = 54; Col = "E"; Value = 164.68 },
    @{ Row = 54; Col = "F"; Value = 174.68 },
    @{ Row = 55; Col = "A"; Value = 46029 },
    @{ Row = 55; Col = "D"; Value = 161.49 },
    @{ Row = 55; Col = "E"; Value = 163.26 },
    @{ Row = 55; Col = "F"; Value = 173.26 },
    @{ Row = 56; Col = "A"; Value = 46029 },
    @{ Row = 56; Col = "D"; Value = 158.56 },
    @{ Row = 57; Col = "A"; Value = 46029 },
    @{ Row = 57; Col = "D"; Value = 159.29 },
    @{ Row = 57; Col = "E"; Value = 157.68 },
    @{ Row = 58; Col = "A"; Value = 46029 },
    @{ Row = 58; Col = "D"; Value = 155.05000000000001 },
    @{ Row = 58; Col = "E"; Value = 153.58000000000001 },
    @{ Row = 58; Col = "F"; Value = 163.58000000000001 },
    @{ Row = 59; Col = "A"; Value = 46029 },
    @{ Row = 59; Col = "D"; Value = 160.91 },
    @{ Row = 59; Col = "E"; Value = 163.18 },
    @{ Row = 60; Col = "A"; Value = 46028 },
    @{ Row = 60; Col = "D"; Value = 168.55 },
    @{ Row = 60; Col = "E"; Value = 164.96 },
    @{ Row = 60; Col = "F"; Value = 174.96 },
    @{ Row = 61; Col = "A"; Value = 46028 },
    @{ Row = 61; Col = "D"; Value = 161.57 },
    @{ Row = 61; Col = "E"; Value = 163.11000000000001 },
    @{ Row = 61; Col = "F"; Value = 173.11 },
    @{ Row = 62; Col = "A"; Value = 46028 },
    @{ Row = 62; Col = "D"; Value = 158.34 },
    @{ Row = 63; Col = "A"; Value = 46028 },
    @{ Row = 63; Col = "D"; Value = 159.1 },
    @{ Row = 63; Col = "E"; Value = 157.53 },
    @{ Row = 64; Col = "A"; Value = 46028 },
    @{ Row = 64; Col = "D"; Value = 154.87 },
    @{ Row = 64; Col = "E"; Value = 153.43 },
    @{ Row = 64; Col = "F"; Value = 163.43 },
    @{ Row = 65; Col = "A"; Value = 46028 },
    @{ Row = 65; Col = "D"; Value = 161.06 },
    @{ Row = 65; Col = "E"; Value = 163.5 }
)

foreach ($e in $edits) {
    $cellRef = "$($e.Col)$($e.Row)"
    $ws.Range($cellRef).Value = $e.Value
}
